# Apply the changes described by the diff:
#  - polyFeature (sheet1): update row 18 (PC) values, add new row 34 (PPSU)
#  - polyIndex (sheet2): selection/active-sheet bookkeeping changes
#  - make polyFeature the active/selected sheet (tabSelected) instead of polyIndex

$wb = $excel.ActiveWorkbook

$wsFeature = $wb.Worksheets.Item("polyFeature")
$wsIndex   = $wb.Worksheets.Item("polyIndex")

# --- Update existing row 18 (PC) on polyFeature ---
$wsFeature.Range("C18").Value = 3
$wsFeature.Range("F18").Value = -3
$wsFeature.Range("K18").Value = 5

# --- Append a new row 34 (PPSU) on polyFeature ---
$wsFeature.Range("A34").Value = "PPSU"
$wsFeature.Range("B34").Value = 0
$wsFeature.Range("C34").Value = 7
$wsFeature.Range("D34").Value = 5
$wsFeature.Range("E34").Value = 0
$wsFeature.Range("F34").Value = -2
$wsFeature.Range("G34").Value = 1
$wsFeature.Range("H34").Value = 0
$wsFeature.Range("I34").Value = 0
$wsFeature.Range("J34").Value = 2
$wsFeature.Range("K34").Value = 0

# Column K has no sheet-wide column style (unlike B:J), so the new K34 cell
# needs its centered style set explicitly to match the rest of the column.
$wsFeature.Range("K34").HorizontalAlignment = -4108

# --- Update selections / view state to match the saved workbook ---
$wsFeature.Range("A20").Select()
$wsFeature.Range("L34").Select()

$wsIndex.Range("B1:I1").Select()

# Make polyFeature the active (selected) sheet/tab
$wsFeature.Activate()
$wsFeature.Select()
